# PROG6001.docx edit
#
# Target change (see commit message "Amritveerkaur1301 updated this document"):
#   1. The paragraph "Ben changing things up!" becomes:
#        - a run "Amritveerkaur1301" carrying GitHub-mention style formatting
#          (Helvetica font, dark-green color, small size)
#        - a run " changing things up!"
#        - a trailing, separate run " "
#   2. A brand-new paragraph "I updated this document" is inserted right
#      after it (before the existing "Again!" paragraph), with the
#      "_GoBack" bookmark sitting between "I upd" and "ated this document".

$d = $word.ActiveDocument

# --- 1. Swap the author's name in place, keeping the rest of the sentence. ---
# Using Find/Replace keeps this robust even if exact run boundaries shift;
# the match is anchored on the whole original sentence so we can't clobber
# anything else in the document.
$d.Content.Find.Execute("Ben changing things up!", $true, $false, $false, $false, $false, $true, 1, $false, "Amritveerkaur1301 changing things up!", 2)

# The sentence lives in paragraph 5 both before and after the replace (no
# paragraphs are added/removed by a Find/Replace), so we can address it
# directly.
$targetPara = $d.Paragraphs(5)

# --- 2. Give just the "Amritveerkaur1301" run its own character formatting. ---
$nameLength = "Amritveerkaur1301".Length
$nameStart = $targetPara.Range.Start
$nameRange = $d.Range($nameStart, $nameStart + $nameLength)
$nameRange.Font.Name = "Helvetica"
$nameRange.Font.Color = 7680
$nameRange.Font.Size = 10.5

# --- 3. Append a standalone trailing-space run at the end of the paragraph. ---
$targetPara = $d.Paragraphs(5)
$beforeMark = $targetPara.Range.End - 1
$tailRange = $d.Range($beforeMark, $beforeMark)
$tailRange.InsertAfter(" ")

# --- 4. Insert the new "I updated this document" paragraph right after it. ---
$targetPara = $d.Paragraphs(5)
$afterPara = $targetPara.Range.End
$newParaRange = $d.Range($afterPara, $afterPara)
$newParaRange.InsertAfter("I updated this document`r")

# --- 5. Drop the "_GoBack" bookmark between "I upd" and "ated this document". ---
$newPara = $d.Paragraphs(6)
$bookmarkPos = $newPara.Range.Start + "I upd".Length
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
